# Top_50_city_data_final.xlsx ETL re-export
#
# The "Pop_Growth" column (H) was recomputed/re-imported with values
# rounded to 2 decimal places, and the ad-hoc number formatting that had
# been applied to the raw population/rent/growth columns (thousands
# separators on D:G, X:Y, AB and the accounting-style format on H) was
# cleared so the sheet falls back to the workbook's default "General"
# format for those columns - matching how the data looked when the
# working database creation files were regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Pop_Growth column (H) to 2 decimal places, in place.
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $v = $cell.Value2
    $cell.Value = [Math]::Round([double]$v, 2)
}

# Clear the inherited number formatting (thousands separators / accounting
# format) from the population, rent, income and growth columns so they
# revert to the default General format.
$ws.Range("D2:G51").ClearFormats()
$ws.Range("H2:H51").ClearFormats()
$ws.Range("X2:Y51").ClearFormats()
$ws.Range("AB2:AB51").ClearFormats()

# Reset the sheet selection back to the top-left cell (it had been left
# selecting the whole Pop_Growth column).
$ws.Range("A1").Select()
